$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Midterm 1 grades (column K) for each student row.
# Cells with a formula (score/50):
$ws.Range("K3").Formula  = "=47/50"
$ws.Range("K4").Formula  = "=31/50"
$ws.Range("K5").Formula  = "=43/50"
$ws.Range("K6").Formula  = "=46.5/50"
$ws.Range("K7").Formula  = "=43/50"
$ws.Range("K11").Formula = "=41/50"
$ws.Range("K12").Formula = "=42/50"
$ws.Range("K14").Formula = "=37.5/50"
$ws.Range("K16").Formula = "=47/50"
$ws.Range("K17").Formula = "=45.5/50"
$ws.Range("K18").Formula = "=42/50"
$ws.Range("K19").Formula = "=45.5/50"

# Rows for students with late/missing submissions - plain 0 value, no formula.
$ws.Range("K8").Value = 0
$ws.Range("K9").Value = 0

# Move the active selection to K7, as recorded after grading.
$ws.Range("K7").Select()
